# Updated symbol list on Tue Dec 13 17:43:58 UTC 2022 with GitHub Actions
#
# Applies the refreshed "Price" (column D) figures and a couple of
# "Volume(1h)" (column E) label tweaks to the crypto ranking sheet.
#
# Column D holds numeric-looking values that are stored as TEXT in the
# workbook (t="inlineStr"/shared-string, not t="n"). Plain `.Value =`
# assignment would let Excel auto-convert a numeric-looking string into a
# real number, which would flip the cell's stored type. Setting
# NumberFormat to the text format ("@") before assigning the value keeps
# these cells text, matching the original authoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
}

# --- Column D: Price updates ---
Set-TextValue "D2"  "270.56"
Set-TextValue "D3"  "22.64"
Set-TextValue "D4"  "6.335"
Set-TextValue "D5"  "0.06189"
Set-TextValue "D6"  "3.659"
Set-TextValue "D7"  "6.658"
Set-TextValue "D8"  "1.389"
Set-TextValue "D9"  "0.8314"
Set-TextValue "D11" "0.1603"
Set-TextValue "D12" "0.08276"
Set-TextValue "D13" "0.03560"
Set-TextValue "D14" "0.03256"
Set-TextValue "D15" "4.068"
Set-TextValue "D16" "0.09307"
Set-TextValue "D17" "0.001657"
Set-TextValue "D18" "0.04752"
Set-TextValue "D19" "0.006366"
Set-TextValue "D20" "0.005669"
Set-TextValue "D21" "0.001079"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "3.728"
Set-TextValue "D25" "0.3334"
Set-TextValue "D27" "0.0002708"
Set-TextValue "D40" "0.04730"
Set-TextValue "D41" "0.006997"
Set-TextValue "D42" "0.1159"
Set-TextValue "D43" "0.003547"
Set-TextValue "D44" "0.01186"
Set-TextValue "D45" "0.00006284"
Set-TextValue "D46" "0.0009915"
Set-TextValue "D48" "0.7834"
Set-TextValue "D49" "0.002318"
Set-TextValue "D51" "0.01242"

# --- Column E: Volume(1h) label tweaks ---
$ws.Range("E8").Value  = "7FTXTokenFTT"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
